# Overwrite old files with RMI version
# - Adds a "New Mexico" label + refreshed date stamp on the About sheet
# - Updates dispatch priority for onshore wind & solar PV from 2 -> 1 on BDPbES
#   (all the year-over-year formula columns recompute automatically)
# - Leaves the BDPbES sheet selected/active with B8 highlighted, matching the
#   author's final view state

$wb = $excel.ActiveWorkbook

$wsAbout  = $wb.Worksheets.Item("About")
$wsBDPbES = $wb.Worksheets.Item("BDPbES")

# --- About sheet: add state label, bump the "as of" date ---------------
$wsAbout.Range("B1").Value = "New Mexico"
$wsAbout.Range("C1").Value = 44463

# --- BDPbES sheet: onshore wind (row 6) & solar PV (row 7) now priority 1
$wsBDPbES.Range("B6").Value = 1
$wsBDPbES.Range("B7").Value = 1

# --- Final selection/active sheet state ---------------------------------
$wsBDPbES.Activate()
$wsBDPbES.Range("B8").Select()
